$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old TestScenario_3 / TestScenario_4 rows (17-20); the Approve/Reject
# block that used to live at rows 16-20 is now rows 12-16 after the inserted steps below.
$ws.Rows("17:20").Delete()

# Rewrite rows 1-16 with the updated TestScenario_1 (extra AnnualRevenue/CustomerPriority__c/
# Website + Submit-for-Approval steps) and TestScenario_2 (Approve/Reject Account) content.
$data = @(
    @('TestScenarioID', 'TestCaseID', 'Description', 'Precondition', 'TestData', 'Steps', 'UserAction', 'ExpectedResult', 'Approved/Rejected', 'ReasonToReject'),
    @('TestScenario_1', 'TestScenario_1.TestCase_1', 'New Account', 'User Needs to Login to Salesforce, from the browser with correct credentials', '', 'Step 1', 'Click Account tab, and click on New button', 'Shows fields to enter to create a new Account', '', ''),
    @('', '', '', '', 'Valid value for required field Name', 'Step 2', 'Enter valid value for Name', 'Value accepted for Name', '', ''),
    @('', '', '', '', 'Valid value for Description', 'Step 3', 'Enter valid value for Description', 'Value accepted for Description', '', ''),
    @('', '', '', '', 'Valid value for AccountNumber', 'Step 4', 'Enter valid value for AccountNumber', 'Value accepted for AccountNumber', '', ''),
    @('', '', '', '', 'Valid value for required field AnnualRevenue, value should be greater or equal 50000.45 to submit for Approval towards to the assigned approver - user:radhapatil@gmail.com', 'Step 5', 'Enter valid value for AnnualRevenue', 'Value accepted for AnnualRevenue', '', ''),
    @('', '', '', '', 'Valid value for CustomerPriority__c', 'Step 6', 'Enter valid value for CustomerPriority__c', 'Value accepted for CustomerPriority__c', '', ''),
    @('', '', '', '', 'Valid value for Website', 'Step 7', 'Enter valid value for Website', 'Value accepted for Website', '', ''),
    @('', '', '', '', '', 'Step 8', 'Click on Save button to save Account with fields', 'Saved changes successfully for the Account', '', ''),
    @('', '', '', '', '', 'Step 9', 'If this record meets the entry criteria then Click on ''Submit for Approval'' button to Submit for Approval.', 'Alert message box will be displayed for confirmation from Salesforce.', '', ''),
    @('', '', '', '', '', 'Step 10', 'Click on ''OK'' button to submit the record for Approval.', 'Unable to Submit for Approval message will be displayed if this record does not meet the entry criteria. Otherwise, this record will be displayed under Approval History section with the status ''Pending''.', '', ''),
    @('TestScenario_2', 'TestScenario_2.TestCase_1', 'Approve/Reject Account', 'User Needs to Login to Salesforce, from the browser with correct credentials', '', 'Step 1', 'Click Account tab, ', '', '', ''),
    @('', '', '', '', '', 'Step 2', 'Click on Account name to Approve/Reject', 'Details of Account successfully displayed', '', ''),
    @('', '', '', '', '', 'Step 3', 'Click on Approve/Reject link from ''Approval History Section'' to Approve/Reject Account', 'Details of ''Approve/Reject Approval Request'' successfully displayed', '', ''),
    @('', '', '', '', '', 'Step 4', 'Enter Comments to notify the user (if any)', 'Value accepted for Comments', '', ''),
    @('', '', '', '', '', 'Step 5', 'Click on Approve/Reject button to Approve/Reject Account', 'Overall status (Approved/Rejected) will be displayed under ''Approval History'' section', '', '')
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

# Column width adjustments (E, F, G, H) to match the wider text now stored there.
# (values chosen so the COM ColumnWidth->stored-width rounding lands as close as
# possible to the target widths 156.410625 / 9.410625 / 93.700625 / 179.840625)
$ws.Columns("E").ColumnWidth = 155.5
$ws.Columns("F").ColumnWidth = 8.5
$ws.Columns("G").ColumnWidth = 92.83333333333333
$ws.Columns("H").ColumnWidth = 179.0
